$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting to the new block A63:G74 (font size 14, thin border, row height 18) ---
# matching the look of the rest of the sheet (same font/border as the preceding data rows)
$newRange = $ws.Range("A63:G74")
$newRange.Font.Size = 14
$newRange.Borders.LineStyle = 1

# Set explicit row height (18pt) to match existing data rows
for ($r = 63; $r -le 74; $r++) {
    $ws.Rows.Item($r).RowHeight = 18
}

# --- Populate month names (ES, PT, IT, CA, FR, RU, EN) ---
$ws.Range("A63").Value = "Enero"
$ws.Range("B63").Value = "Janeiro"
$ws.Range("C63").Value = "Gennaio"
$ws.Range("D63").Value = "Gener"
$ws.Range("E63").Value = "Janvier"
$ws.Range("F63").Value = "Январь"
$ws.Range("G63").Value = "January"
$ws.Range("A64").Value = "Febrero"
$ws.Range("B64").Value = "Fevereiro"
$ws.Range("C64").Value = "Febbraio"
$ws.Range("D64").Value = "Febrer"
$ws.Range("E64").Value = "Février"
$ws.Range("F64").Value = "Февраль"
$ws.Range("G64").Value = "February"
$ws.Range("A65").Value = "Marzo"
$ws.Range("B65").Value = "Março"
$ws.Range("C65").Value = "Marzo"
$ws.Range("D65").Value = "Març"
$ws.Range("E65").Value = "Mars"
$ws.Range("F65").Value = "Март"
$ws.Range("G65").Value = "March"
$ws.Range("A66").Value = "Abril"
$ws.Range("B66").Value = "Abril"
$ws.Range("C66").Value = "Aprile"
$ws.Range("D66").Value = "Abril"
$ws.Range("E66").Value = "Avril"
$ws.Range("F66").Value = "Апрель"
$ws.Range("G66").Value = "April"
$ws.Range("A67").Value = "Mayo"
$ws.Range("B67").Value = "Maio"
$ws.Range("C67").Value = "Maggio"
$ws.Range("D67").Value = "Maig"
$ws.Range("E67").Value = "Mai"
$ws.Range("F67").Value = "Май"
$ws.Range("G67").Value = "May"
$ws.Range("A68").Value = "Junio"
$ws.Range("B68").Value = "Junho"
$ws.Range("C68").Value = "Giugno"
$ws.Range("D68").Value = "Juny"
$ws.Range("E68").Value = "Juin"
$ws.Range("F68").Value = "Июнь"
$ws.Range("G68").Value = "June"
$ws.Range("A69").Value = "Julio"
$ws.Range("B69").Value = "Julho"
$ws.Range("C69").Value = "Iuglio"
$ws.Range("D69").Value = "Juliol"
$ws.Range("E69").Value = "Juillet"
$ws.Range("F69").Value = "Июль"
$ws.Range("G69").Value = "July"
$ws.Range("A70").Value = "Agosto"
$ws.Range("B70").Value = "Agosto"
$ws.Range("C70").Value = "Agosto"
$ws.Range("D70").Value = "Agost"
$ws.Range("E70").Value = "Août"
$ws.Range("F70").Value = "Август"
$ws.Range("G70").Value = "August"
$ws.Range("A71").Value = "Septiembre"
$ws.Range("B71").Value = "Setembre"
$ws.Range("C71").Value = "Settembre"
$ws.Range("D71").Value = "Setembre"
$ws.Range("E71").Value = "Septembre"
$ws.Range("F71").Value = "Сентябрь"
$ws.Range("G71").Value = "September"
$ws.Range("A72").Value = "Octubre"
$ws.Range("B72").Value = "Outubro"
$ws.Range("C72").Value = "Ottobre"
$ws.Range("D72").Value = "Octubre"
$ws.Range("E72").Value = "Octubre"
$ws.Range("F72").Value = "Октябрь"
$ws.Range("G72").Value = "October"
$ws.Range("A73").Value = "Noviembre"
$ws.Range("B73").Value = "Novembro"
$ws.Range("C73").Value = "Novembre"
$ws.Range("D73").Value = "Novembre"
$ws.Range("E73").Value = "Novembre"
$ws.Range("F73").Value = "Ноябрь"
$ws.Range("G73").Value = "November"
$ws.Range("A74").Value = "Diciembre"
$ws.Range("B74").Value = "Dezembro"
$ws.Range("C74").Value = "Dicembre"
$ws.Range("D74").Value = "Desembre"
$ws.Range("E74").Value = "Décembre"
$ws.Range("F74").Value = "Декабрь"
$ws.Range("G74").Value = "December"

# --- Update sheet view / selection to match the post-edit state ---
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F64").Select()

